$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data rows (2-4) entirely, shifting rows up
$ws.Range("A2:A4").EntireRow.Delete()

# Remove the now-unneeded columns C:E entirely, shifting columns left
$ws.Range("C1:E1").EntireColumn.Delete()

# Update the remaining header cell (B1) with the new label
$ws.Range("B1").Value = "Tafel slope [V/dec]"
